$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 400 (rows 400..462 shift down to 401..463)
$ws.Rows("400:400").Insert()

# Populate the new row 400 with the new data record
$r = 400
$ws.Cells.Item($r, 1).Value = 4
$ws.Cells.Item($r, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item($r, 3).Value = 'Los Lagos'
$ws.Cells.Item($r, 4).Value = 45180
$ws.Cells.Item($r, 5).Value = 10
$ws.Cells.Item($r, 6).Value = 100112032
$ws.Cells.Item($r, 7).Value = 'Zapallo italiano'
$ws.Cells.Item($r, 8).Value = 'Sin especificar'
$ws.Cells.Item($r, 9).Value = 'Primera'
$ws.Cells.Item($r, 10).Value = 70
$ws.Cells.Item($r, 11).Value = 20000
$ws.Cells.Item($r, 12).Value = 20000
$ws.Cells.Item($r, 13).Value = 20000
$ws.Cells.Item($r, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item($r, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item($r, 16).Value = 400
$ws.Cells.Item($r, 17).Value = 50
$ws.Cells.Item($r, 18).Value = 'Hortaliza'
